$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66; existing rows 66-139 shift down to 67-140.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new data record.
$ws.Range("A66").Value = 7
$ws.Range("B66").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C66").Value = "Ñuble"
$ws.Range("D66").Value = 44664
$ws.Range("E66").Value = 16
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100109
$ws.Range("H66").Value = "Uva"
$ws.Range("I66").Value = 100109001
$ws.Range("J66").Value = "Uva"
$ws.Range("K66").Value = "Thompson seedless"
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 120
$ws.Range("N66").Value = 11000
$ws.Range("O66").Value = 12000
$ws.Range("P66").Value = 11500
$ws.Range("Q66").Value = "$/bandeja 18 kilos"
$ws.Range("R66").Value = "Región de O'Higgins"
$ws.Range("S66").Value = 639
$ws.Range("T66").Value = 18
